$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.947566
$ws.Range("H2").Value = 2.842698
$ws.Range("I2").Value = 0.01860557227824198
$ws.Range("J2").Value = 0.01860557227824198
$ws.Range("M2").Value = 2.759544333333333
$ws.Range("N2").Value = 8.278632999999999
$ws.Range("O2").Value = 0.2574067337278401
$ws.Range("P2").Value = 0.2574067337278401
$ws.Range("Q2").Value = 2.614850385759333
$ws.Range("R2").Value = 23.533653471834
$ws.Range("S2").Value = 0.004789199589279517
$ws.Range("T2").Value = 0.004789199589279516
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.947566
$ws.Range("H3").Value = 2.842698
$ws.Range("I3").Value = 0.01860557227824198
$ws.Range("J3").Value = 0.01860557227824198
$ws.Range("O3").Value = 0.6758254232987829
$ws.Range("P3").Value = 0.6758254232987829
$ws.Range("Q3").Value = 6.865330767482
$ws.Range("R3").Value = 61.787976907338
$ws.Range("S3").Value = 0.01257411876065899
$ws.Range("T3").Value = 0.01257411876065899
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.947566
$ws.Range("H4").Value = 2.842698
$ws.Range("I4").Value = 0.01860557227824198
$ws.Range("J4").Value = 0.01860557227824198
$ws.Range("M4").Value = 0.5200313333333334
$ws.Range("N4").Value = 1.560094
$ws.Range("O4").Value = 0.0485078515798926
$ws.Range("P4").Value = 0.0485078515798926
$ws.Range("Q4").Value = 0.4927640104013334
$ws.Range("R4").Value = 4.434876093612
$ws.Range("S4").Value = 0.0009025163386319263
$ws.Range("T4").Value = 0.0009025163386319261
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.947566
$ws.Range("H5").Value = 2.842698
$ws.Range("I5").Value = 0.01860557227824198
$ws.Range("J5").Value = 0.01860557227824198
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1957573333333333
$ws.Range("N5").Value = 0.587272
$ws.Range("O5").Value = 0.01825999139348442
$ws.Range("P5").Value = 0.01825999139348442
$ws.Range("Q5").Value = 0.1854929933173333
$ws.Range("R5").Value = 1.669436939856
$ws.Range("S5").Value = 0.000339737589671551
$ws.Range("T5").Value = 0.0003397375896715509
$ws.Range("I6").Value = 0.9615746457924133
$ws.Range("J6").Value = 0.9615746457924131
$ws.Range("M6").Value = 2.759544333333333
$ws.Range("N6").Value = 8.278632999999999
$ws.Range("O6").Value = 0.2574067337278401
$ws.Range("P6").Value = 0.2574067337278401
$ws.Range("Q6").Value = 135.1409027298281
$ws.Range("R6").Value = 1216.268124568453
$ws.Range("S6").Value = 0.2475157888089299
$ws.Range("T6").Value = 0.2475157888089298
$ws.Range("I7").Value = 0.9615746457924133
$ws.Range("J7").Value = 0.9615746457924131
$ws.Range("O7").Value = 0.6758254232987829
$ws.Range("P7").Value = 0.6758254232987829
$ws.Range("S7").Value = 0.649856592026035
$ws.Range("T7").Value = 0.6498565920260349
$ws.Range("I8").Value = 0.9615746457924133
$ws.Range("J8").Value = 0.9615746457924131
$ws.Range("M8").Value = 0.5200313333333334
$ws.Range("N8").Value = 1.560094
$ws.Range("O8").Value = 0.0485078515798926
$ws.Range("P8").Value = 0.0485078515798926
$ws.Range("Q8").Value = 25.46706823498378
$ws.Range("R8").Value = 229.203614114854
$ws.Range("S8").Value = 0.04664392020108618
$ws.Range("T8").Value = 0.04664392020108617
$ws.Range("I9").Value = 0.9615746457924133
$ws.Range("J9").Value = 0.9615746457924131
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.1957573333333333
$ws.Range("N9").Value = 0.587272
$ws.Range("O9").Value = 0.01825999139348442
$ws.Range("P9").Value = 0.01825999139348442
$ws.Range("Q9").Value = 9.586663429572445
$ws.Range("R9").Value = 86.279970866152
$ws.Range("S9").Value = 0.0175583447563623
$ws.Range("T9").Value = 0.0175583447563623
$ws.Range("G10").Value = 0.8226676666666667
$ws.Range("H10").Value = 2.468003
$ws.Range("I10").Value = 0.01615317849431
$ws.Range("J10").Value = 0.01615317849431
$ws.Range("M10").Value = 2.759544333333333
$ws.Range("N10").Value = 8.278632999999999
$ws.Range("O10").Value = 0.2574067337278401
$ws.Range("P10").Value = 0.2574067337278401
$ws.Range("Q10").Value = 2.270187897766555
$ws.Range("R10").Value = 20.431691079899
$ws.Range("S10").Value = 0.004157936915543127
$ws.Range("T10").Value = 0.004157936915543126
$ws.Range("G11").Value = 0.8226676666666667
$ws.Range("H11").Value = 2.468003
$ws.Range("I11").Value = 0.01615317849431
$ws.Range("J11").Value = 0.01615317849431
$ws.Range("O11").Value = 0.6758254232987829
$ws.Range("P11").Value = 0.6758254232987829
$ws.Range("Q11").Value = 5.960413990560333
$ws.Range("R11").Value = 53.643725915043
$ws.Range("S11").Value = 0.01091672869353785
$ws.Range("T11").Value = 0.01091672869353785
$ws.Range("G12").Value = 0.8226676666666667
$ws.Range("H12").Value = 2.468003
$ws.Range("I12").Value = 0.01615317849431
$ws.Range("J12").Value = 0.01615317849431
$ws.Range("M12").Value = 0.5200313333333334
$ws.Range("N12").Value = 1.560094
$ws.Range("O12").Value = 0.0485078515798926
$ws.Range("P12").Value = 0.0485078515798926
$ws.Range("Q12").Value = 0.427812963586889
$ws.Range("R12").Value = 3.850316672282
$ws.Range("S12").Value = 0.0007835559849455026
$ws.Range("T12").Value = 0.0007835559849455022
$ws.Range("G13").Value = 0.8226676666666667
$ws.Range("H13").Value = 2.468003
$ws.Range("I13").Value = 0.01615317849431
$ws.Range("J13").Value = 0.01615317849431
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.1957573333333333
$ws.Range("N13").Value = 0.587272
$ws.Range("O13").Value = 0.01825999139348442
$ws.Range("P13").Value = 0.01825999139348442
$ws.Range("Q13").Value = 0.1610432286462222
$ws.Range("R13").Value = 1.449389057816
$ws.Range("S13").Value = 0.0002949569002835182
$ws.Range("T13").Value = 0.0002949569002835182
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.186737
$ws.Range("H14").Value = 0.560211
$ws.Range("I14").Value = 0.003666603435034682
$ws.Range("J14").Value = 0.003666603435034681
$ws.Range("M14").Value = 2.759544333333333
$ws.Range("N14").Value = 8.278632999999999
$ws.Range("O14").Value = 0.2574067337278401
$ws.Range("P14").Value = 0.2574067337278401
$ws.Range("Q14").Value = 0.5153090301736666
$ws.Range("R14").Value = 4.637781271563
$ws.Range("S14").Value = 0.0009438084140875561
$ws.Range("T14").Value = 0.0009438084140875559
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.186737
$ws.Range("H15").Value = 0.560211
$ws.Range("I15").Value = 0.003666603435034682
$ws.Range("J15").Value = 0.003666603435034681
$ws.Range("O15").Value = 0.6758254232987829
$ws.Range("P15").Value = 0.6758254232987829
$ws.Range("Q15").Value = 1.352951954299
$ws.Range("R15").Value = 12.176567588691
$ws.Range("S15").Value = 0.002477983818551085
$ws.Range("T15").Value = 0.002477983818551085
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.186737
$ws.Range("H16").Value = 0.560211
$ws.Range("I16").Value = 0.003666603435034682
$ws.Range("J16").Value = 0.003666603435034681
$ws.Range("M16").Value = 0.5200313333333334
$ws.Range("N16").Value = 1.560094
$ws.Range("O16").Value = 0.0485078515798926
$ws.Range("P16").Value = 0.0485078515798926
$ws.Range("Q16").Value = 0.09710909109266669
$ws.Range("R16").Value = 0.873981819834
$ws.Range("S16").Value = 0.0001778590552289867
$ws.Range("T16").Value = 0.0001778590552289867
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.186737
$ws.Range("H17").Value = 0.560211
$ws.Range("I17").Value = 0.003666603435034682
$ws.Range("J17").Value = 0.003666603435034681
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.1957573333333333
$ws.Range("N17").Value = 0.587272
$ws.Range("O17").Value = 0.01825999139348442
$ws.Range("P17").Value = 0.01825999139348442
$ws.Range("Q17").Value = 0.03655513715466667
$ws.Range("R17").Value = 0.328996234392
$ws.Range("S17").Value = 0.00006695214716705371
$ws.Range("T17").Value = 0.0000669521471670537
